$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row 2 (Id=2, "В дневном стационаре"): Options column F2 changes from NULL to |DayStationary|
$ws.Range("F2").Value = "|DayStationary|"

# Row 5 (Id=6, "Стационарно"): Options column F5 changes from NULL to |Stationary|
$ws.Range("F5").Value = "|Stationary|"
